# Generate Report for Handoff
# - Priority for the 4 "Ready for handoff" rows moves from "low" to "ht"
#   (both the zh-cn and de-de localization-status sheets).
# - The zh-cn sheet's "Latest Handoff Datetime" for those same rows is
#   refreshed to the new handoff timestamp.
# - The de-de sheet (and the Overview rollup) shared the older handoff
#   timestamp string, which is refreshed the same way so every cell that
#   displayed it picks up the new value.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")
$ov = $wb.Worksheets.Item("Overview")

# Priority: low -> ht, rows 4-7 (339a660f.., a36dfb4e.., b4e77864.., e82f82ba..)
$zh.Range("E4:E7").Value = "ht"
$de.Range("E4:E7").Value = "ht"

# zh-cn Latest Handoff Datetime refresh for the same rows
$zh.Range("H4:H7").Value = "2016-08-12 02:53:30"

# The de-de sheet and Overview rollup still show the older handoff
# generation timestamp for those rows; refresh it to the new value.
$de.Range("H4:H7").Value = "2016-08-12 02:53:36"
$ov.Range("G4:G7").Value = "2016-08-12 02:53:36"
